$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: new "Anmerkung" column (G) ---
$ws.Range("G1").Value = "Anmerkung"

# --- Row 2: taller row for wrapped header-ish content already there ---
$ws.Rows.Item(2).RowHeight = 54.75

# --- Row 5: Brunel row gets a status note ---
$ws.Range("E5").Value = "17.09 Tel-Interview"
$ws.Range("F5").Value = "Gespräch am 22.09"

# --- Row 10: MBTech row gets a status note ---
$ws.Range("E10").Value = "10.09 Inteview"
$ws.Range("F10").Value = "Im Lauf"

# --- Row 11: Bosch row gets a status note ---
$ws.Range("E11").Value = 18.09
$ws.Range("F11").Value = "Abgesagt"

# --- New row 13: Siemens Beijing ---
$ws.Range("A13").Value = "Siemens"
$ws.Range("B13").Value = "Siemens Graduate Program – Digital Factory Division – Business Development –"
$ws.Range("C13").Value = " Beijing"
$ws.Range("D13").Value = "19.09.15"
$ws.Range("G13").Value = "yingjiesheng.com-han_tiger7"
$ws.Rows.Item(13).RowHeight = 29.25

$r = $ws.Range("B13")
$r.Font.Name = "Verdana"
$r.Font.Size = 9
$r.Font.Color = 0
$r.WrapText = $true

# --- New row 14: Siemens Shanghai ---
$ws.Range("A14").Value = "Siemens"
$ws.Range("B14").Value = "Siemens Graduate Program–Siemens Wind Power Blades –Production Management-Shanghai"
$ws.Range("C14").Value = "Shanghai"
$ws.Range("D14").Value = "19.09.15"
$ws.Rows.Item(14).RowHeight = 22.5

$r = $ws.Range("B14")
$r.Font.Name = "Verdana"
$r.Font.Size = 9
$r.Font.Color = 0
$r.WrapText = $true

# --- New row 15: Bosch Shanghai ---
$ws.Range("A15").Value = "Bosch"
$ws.Range("B15").Value = "Algorithm Development"
$ws.Range("C15").Value = "Shanghai"
$ws.Range("D15").Value = "19.09.15"
$ws.Range("G15").Value = "51job han_tiger7"

# --- New row 16: job-ad footer note ---
$r = $ws.Range("B16")
$r.Value = " Software Engineer (Radar) 软件工程师 （雷达）(汽车底盘控制 )"
$r.Font.Name = "Arial, Helvetica, sans-serif;宋体"
$r.Font.Size = 9
$r.Font.Color = 7423762

# --- Column widths ---
$ws.Columns.Item(4).ColumnWidth = 20.714285714285715
$ws.Columns.Item(5).ColumnWidth = 21.142857142857146
$ws.Columns.Item(6).ColumnWidth = 20.57142857142857
$ws.Columns.Item(7).ColumnWidth = 20.42857142857143

# --- Selection / view state ---
$ws.Range("C17").Select()

Write-Output "edit applied"
